$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$csv = @"
1202,301,V,Dallas,7,6,7,7,27,4.5,7,-210,3
1202,302,H,NewOrleans,0,7,3,7,17,48,45.5,180,22.5
1205,451,V,Minnesota,6,0,9,12,27,7.5,7.5,-310,7
1205,452,H,Detroit,0,20,3,6,29,48,47.5,260,25
1205,453,V,Arizona,14,7,3,9,33,7,9,-370,4
1205,454,H,Chicago,0,7,7,8,22,46.5,43,310,20
1205,455,V,TampaBay,13,7,7,3,30,9.5,12,-550,4
1205,456,H,Atlanta,7,10,0,0,17,50.5,50.5,425,26
1205,457,V,Denver,0,3,0,6,9,50.5,46.5,300,23
1205,458,H,KansasCity,10,0,3,9,22,9.5,9.5,-360,4
1205,459,V,Indianapolis,7,7,7,10,31,7,11,-500,5
1205,460,H,Houston,0,0,0,0,0,47,45.5,400,21
1205,461,V,Philadelphia,7,17,3,6,33,6.5,4.5,-210,1.5
1205,462,H,NYJets,12,6,0,0,18,45,44,180,23.5
1205,463,V,LAChargers,16,8,0,17,41,48,49.5,125,26.5
1205,464,H,Cincinnati,0,13,9,0,22,1,3,-145,4.5
1205,465,V,NYGiants,0,3,3,3,9,43.5,40,260,19.5
1205,466,H,Miami,0,10,0,10,20,2.5,7.5,-310,3
1205,467,V,Washington,7,0,0,10,17,48,47.5,115,22.5
1205,468,H,LasVegas,0,3,3,9,15,1.5,1,-135,1.5
1205,469,V,Baltimore,0,7,3,9,19,2.5,3,-190,1.5
1205,470,H,Pittsburgh,0,3,0,17,20,45,44,170,21
1205,471,V,Jacksonville,0,7,0,0,7,48,48,600,21.5
1205,472,H,LARams,10,6,14,7,37,12,14.5,-900,5.5
1205,473,V,SanFrancisco,14,9,0,0,23,46,3,-145,2
1205,474,H,Seattle,7,14,9,0,30,1,45,125,22.5
1206,475,V,NewEngland,8,3,0,3,14,46,39.5,130,18
1206,476,H,Buffalo,7,0,3,0,10,3.5,3,-150,1.5
"@

$lines = $csv -split "`n" | Where-Object { $_.Trim() -ne "" }

$startRow = 362
$rowIndex = $startRow

foreach ($line in $lines) {
    $fields = $line.Trim() -split ","

    $ws.Cells.Item($rowIndex, 1).Value = [double]$fields[0]
    $ws.Cells.Item($rowIndex, 2).Value = [double]$fields[1]
    $ws.Cells.Item($rowIndex, 3).Value = $fields[2]
    $ws.Cells.Item($rowIndex, 4).Value = $fields[3]
    $ws.Cells.Item($rowIndex, 5).Value = [double]$fields[4]
    $ws.Cells.Item($rowIndex, 6).Value = [double]$fields[5]
    $ws.Cells.Item($rowIndex, 7).Value = [double]$fields[6]
    $ws.Cells.Item($rowIndex, 8).Value = [double]$fields[7]
    $ws.Cells.Item($rowIndex, 9).Value = [double]$fields[8]
    $ws.Cells.Item($rowIndex, 10).Value = [double]$fields[9]
    $ws.Cells.Item($rowIndex, 11).Value = [double]$fields[10]
    $ws.Cells.Item($rowIndex, 12).Value = [double]$fields[11]
    $ws.Cells.Item($rowIndex, 13).Value = [double]$fields[12]

    $rowIndex++
}
